$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.340.20'
$ws.Range("E2").Value = '  -2.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.49'
$ws.Range("E3").Value = '  -3.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.12'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4732'
$ws.Range("E7").Value = '  +1.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3957'
$ws.Range("E8").Value = '  -1.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.89'
$ws.Range("E9").Value = '  -11.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07996'
$ws.Range("E10").Value = '  -4.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.014'
$ws.Range("E11").Value = '  -3.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.50'
$ws.Range("E12").Value = '  -3.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.865.46'
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.946'
$ws.Range("E14").Value = '  -2.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.140'
$ws.Range("E15").Value = '  -4.04%  '
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.26'
$ws.Range("E17").Value = '  -3.81%  '
$ws.Range("E18").Value = '  -3.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06554'
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.26'
$ws.Range("E20").Value = '  -3.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.485'
$ws.Range("E22").Value = '  -4.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.352.10'
$ws.Range("E23").Value = '  -2.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.92'
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.296'
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.072.97'
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.37'
$ws.Range("E27").Value = '  +1.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.97'
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.073'
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.489'
$ws.Range("E30").Value = '  -4.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.11'
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09513'
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9520'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.444'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.588'
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.276'
$ws.Range("E36").Value = '  -4.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06053'
$ws.Range("E37").Value = '  -2.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02232'
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.211'
$ws.Range("E39").Value = '  -4.19%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.012'
$ws.Range("E41").Value = '  -9.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5936'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1897'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.27'
$ws.Range("E44").Value = '  -7.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.279'
$ws.Range("E45").Value = '  -2.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5650'
$ws.Range("E46").Value = '  -3.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.16'
$ws.Range("E47").Value = '  -4.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.428'
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.925'
$ws.Range("E49").Value = '  -4.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06762'
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.42'
$ws.Range("E51").Value = '  -1.86%  '
